$d = $word.ActiveDocument

# Locate the paragraph that currently holds the (hidden) "_GoBack" bookmark -
# this is the empty paragraph right after "Occam's Razor" that we need to
# turn into the new "Model Options" / "Logistic Regression..." content.
$bm = $d.Bookmarks.Item("_GoBack")
$target = $bm.Range.Paragraphs.Item(1)
$targetIndex = $target.Index
$targetRange = $target.Range

# Split it into three paragraphs, keeping the bookmark on the trailing one:
#   1) "Model Options"                                   (bold)
#   2) "Logistic Regression: find decision boundary between classes"
#   3) (empty, holds the _GoBack bookmark - unchanged)
$targetRange.InsertParagraphBefore()
$targetRange.InsertParagraphBefore()

$modelOptionsPara = $d.Paragraphs.Item($targetIndex)
$logisticPara = $d.Paragraphs.Item($targetIndex + 1)

$modelOptionsPara.Range.Text = "Model Options"
$logisticPara.Range.Text = "Logistic Regression: find decision boundary between classes"

# Bold the new "Model Options" heading line.
$modelOptionsPara.Range.Bold = 1
$modelOptionsPara.Range.BoldBi = 1

Write-Host "Inserted Model Options / Logistic Regression paragraphs before the bookmark."
